$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Summary table (rows 10-12): recompute Right/Wrong/Not-Attempt/Max and
#    totals now that the student's (float-safe) answers have been graded.
# ---------------------------------------------------------------------------

# Give the row-label cells (A10/A11/A12) the same "mtitleStyle" look already
# used by the other section headers (A9, A15 ...), without touching their text.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B10").Value = 18
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 72
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "67/112"

# ---------------------------------------------------------------------------
# 2) Per-question "Student Ans" columns.
#    Column A mirrors question-set 1 (graded against column B).
#    Column D keeps only the first three rows of question-set 2 (graded
#    against column E); question-set 3 (columns G/H) is dropped entirely.
# ---------------------------------------------------------------------------

# Style anchor cells already present in the sheet.
$correctAnchor = "B10"   # correctStyle   (green)
$incorrectAnchor = "C10" # incorrectStyle (red)

# row -> student answer for column A ("" means left blank / not attempted)
$colA = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C";
    20 = "Option B"; 21 = "Option C"; 22 = "";          23 = "";
    24 = "Option B"; 25 = "Option A"; 26 = "Option C"; 27 = "";
    28 = "Option B"; 29 = "Option D"; 30 = "";          31 = "Option D";
    32 = "Option C"; 33 = "Option D"; 34 = "Option B"; 35 = "";
    36 = "Option D"; 37 = "Option A"; 38 = "Option A"; 39 = "Option D";
    40 = "Option B"
}

for ($row = 16; $row -le 40; $row++) {
    $answer = $colA[$row]
    $correct = $ws.Range("B$row").Value2

    if ([string]::IsNullOrEmpty($answer)) {
        # Leave blank; cell already carries the default "normalStyle" look.
        continue
    }

    if ($answer -eq $correct) {
        $ws.Range($correctAnchor).Copy() | Out-Null
    } else {
        $ws.Range($incorrectAnchor).Copy() | Out-Null
    }
    $ws.Range("A$row").PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("A$row").Value = $answer
}

# Column D: only rows 16-18 keep an answer now.
$ws.Range($correctAnchor).Copy() | Out-Null
$ws.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D16").Value = "Option A"

$ws.Range($incorrectAnchor).Copy() | Out-Null
$ws.Range("D17").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D17").Value = "Option B"

$ws.Range($correctAnchor).Copy() | Out-Null
$ws.Range("D18").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D18").Value = "Option D"

# ---------------------------------------------------------------------------
# 3) Drop question-set 2's remaining rows (D19:E40) and question-set 3
#    entirely (G15:H40) - these columns are no longer part of the report.
# ---------------------------------------------------------------------------
$ws.Range("D19:E40").Clear() | Out-Null
$ws.Range("G15:H40").Clear() | Out-Null
